$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 385  # was 384
$ws.Cells.Item(4, 6).Value = 152  # was 151
$ws.Cells.Item(5, 6).Value = 1310  # was 1308
$ws.Cells.Item(6, 6).Value = 227  # was 226
$ws.Cells.Item(7, 6).Value = 2496  # was 2492
$ws.Cells.Item(8, 6).Value = 903  # was 897
$ws.Cells.Item(9, 6).Value = 18646  # was 18627
$ws.Cells.Item(10, 6).Value = 52  # was 51
$ws.Cells.Item(11, 6).Value = 1910  # was 1906
$ws.Cells.Item(12, 6).Value = 664  # was 661
$ws.Cells.Item(14, 6).Value = 331  # was 330
$ws.Cells.Item(15, 6).Value = 603  # was 602
$ws.Cells.Item(16, 6).Value = 198  # was 197
$ws.Cells.Item(18, 6).Value = 71  # was 70
$ws.Cells.Item(20, 6).Value = 162  # was 161

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(8, 6).Value = 2  # was 1
$ws.Cells.Item(10, 6).Value = 225  # was 123
$ws.Cells.Item(16, 6).Value = 69  # was 68

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 5888  # was 5885
$ws.Cells.Item(3, 6).Value = 570  # was 569

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(3, 6).Value = 5888  # was 5885
$ws.Cells.Item(4, 6).Value = 570  # was 569
$ws.Cells.Item(6, 6).Value = 385  # was 384
$ws.Cells.Item(8, 6).Value = 152  # was 151
$ws.Cells.Item(10, 6).Value = 1310  # was 1308
$ws.Cells.Item(12, 6).Value = 227  # was 226
$ws.Cells.Item(15, 6).Value = 2496  # was 2492
$ws.Cells.Item(16, 6).Value = 903  # was 897
$ws.Cells.Item(17, 6).Value = 18646  # was 18627
$ws.Cells.Item(19, 6).Value = 2  # was 1
$ws.Cells.Item(20, 6).Value = 52  # was 51
$ws.Cells.Item(22, 6).Value = 225  # was 123
$ws.Cells.Item(23, 6).Value = 225  # was 123
$ws.Cells.Item(24, 6).Value = 1911  # was 1906
$ws.Cells.Item(25, 6).Value = 664  # was 661
$ws.Cells.Item(28, 6).Value = 331  # was 330
$ws.Cells.Item(29, 6).Value = 603  # was 602
$ws.Cells.Item(30, 6).Value = 198  # was 197
$ws.Cells.Item(33, 6).Value = 71  # was 70
$ws.Cells.Item(38, 6).Value = 69  # was 68
$ws.Cells.Item(39, 6).Value = 162  # was 161
